# Re-ran resolve and classify+summarise steps after changes to mapping file.
# This updates the summary tables that are derived from per-species mapping
# data: with the (changed) mapping file, several summaries collapse to zero
# / shift rows because the underlying species sets changed.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "Range Status" sheet: range-size bucket counts all reset to 0 and the
# percentage column (no species counted) is cleared out, row by row.
# ------------------------------------------------------------------
$wsRange = $wb.Worksheets.Item("Range Status")
$wsRange.Range("B2").Value = 0
$wsRange.Range("C2").ClearContents()
$wsRange.Range("B3").Value = 0
$wsRange.Range("C3").ClearContents()
$wsRange.Range("B4").Value = 0
$wsRange.Range("C4").ClearContents()
$wsRange.Range("B5").Value = 0
$wsRange.Range("C5").ClearContents()
$wsRange.Range("B6").Value = 0
$wsRange.Range("C6").ClearContents()
$wsRange.Range("B7").Value = 0
$wsRange.Range("C7").ClearContents()

# ------------------------------------------------------------------
# "Species qualification" sheet: Range Analysis species count drops to 0.
# ------------------------------------------------------------------
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# ------------------------------------------------------------------
# "High Priority break-up" sheet: the "Trend Different" row disappears
# (species re-classified), shifting "IUCN" up, and the remaining rows'
# counts are recomputed.
# ------------------------------------------------------------------
$wsBreak = $wb.Worksheets.Item("High Priority break-up")
$wsBreak.Rows.Item(3).Delete()

$wsBreak.Range("B2").Value = 25
$wsBreak.Range("C2").Value = 55.6
$wsBreak.Range("D2").Value = 25
$wsBreak.Range("E2").Value = 55.6

$wsBreak.Range("A3").Value = "IUCN"
$wsBreak.Range("B3").Value = 20
$wsBreak.Range("C3").Value = 44.4
$wsBreak.Range("D3").Value = 20
$wsBreak.Range("E3").Value = 44.4
